$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the demo UUID values (productID/brand/category/type/unit ids) ---
# Replace the old placeholder GUID values with freshly generated ones,
# keeping the same cell layout / order as before.
$ws.Range("D2").Value = "f243adf4-1cf7-40e7-941c-9549dc987b95"
$ws.Range("I2").Value = "f6e4a2b9-bdd7-44de-97b9-b5f6f2588a20"
$ws.Range("J2").Value = "cbf7c588-eb3f-4a12-9f03-0c71d4365c25"
$ws.Range("K2").Value = "959db27b-15be-471c-acaa-e0ea81168803"
$ws.Range("L2").Value = "7e6d3758-aad6-4e5c-948b-9d15303d0e1c"

# --- Normalize the "Aptos Narrow" font used on column A/D example cells ---
# These cells used a duplicate font entry; reapply the default font so it
# reuses the existing (deduplicated) font definition instead.
$sample = $ws.Range("E2")
$ws.Range("A2:A7").Style = $sample.Style
$ws.Range("D2:D7").Style = $sample.Style

# --- Update the active selection to D3 ---
$ws.Range("D3").Select()
